$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append one more blank "filler" row at the bottom of the sheet (row 444), mirroring
# the L:M-only rows used as spacers further up (rows 242-443 before this edit).
$ws.Rows.Item(444).Insert()
$ws.Cells.Item(444, 12).NumberFormat = "0.00"
$ws.Cells.Item(444, 13).NumberFormat = "0.00"

# --- Clear the stray cells that existed only because of the old row layout; they have
# no counterpart once the table is re-sorted with the two new polls (3/4 and 3/7).
$ws.Cells.Item(241, 21).ClearContents()
$ws.Cells.Item(241, 34).ClearContents()

# --- Rewrite rows 239-246: three pre-existing polls (136, 137, 138) whose stats were
# refreshed by the new rolling-average calculation, plus the five freshly entered polls
# (139-143) that slot in by date among them.

# Row 239 -- poll id 136
$ws.Cells.Item(239, 1).Value = 136
$ws.Cells.Item(239, 2).Value = 2022
$ws.Cells.Item(239, 3).Value = 2
$ws.Cells.Item(239, 4).Value = 28
$ws.Cells.Item(239, 5).Value = 3
$ws.Cells.Item(239, 6).Value = 3
$ws.Cells.Item(239, 7).Value = "opinionway"
$ws.Cells.Item(239, 8).Value = "partially"
$ws.Cells.Item(239, 9).Value = "rolling"
$ws.Cells.Item(239, 10).Value = 1037
$ws.Cells.Item(239, 11).Value = 0
$ws.Cells.Item(239, 12).Value = 1
$ws.Cells.Item(239, 13).Value = 1
$ws.Cells.Item(239, 14).Value = 1
$ws.Cells.Item(239, 15).Value = 1
$ws.Cells.Item(239, 16).Value = 10
$ws.Cells.Item(239, 17).Value = 4
$ws.Cells.Item(239, 19).Value = 5
$ws.Cells.Item(239, 20).Value = 2
$ws.Cells.Item(239, 21).Value = 2
$ws.Cells.Item(239, 22).Value = 27
$ws.Cells.Item(239, 23).Value = 14
$ws.Cells.Item(239, 26).Value = 2
$ws.Cells.Item(239, 27).Value = 2
$ws.Cells.Item(239, 28).Value = 18
$ws.Cells.Item(239, 29).Value = 12
$ws.Cells.Item(239, 12).NumberFormat = "0.00"
$ws.Cells.Item(239, 13).NumberFormat = "0.00"

# Row 240 -- poll id 137
$ws.Cells.Item(240, 1).Value = 137
$ws.Cells.Item(240, 2).Value = 2022
$ws.Cells.Item(240, 3).Value = 2
$ws.Cells.Item(240, 4).Value = 28
$ws.Cells.Item(240, 5).Value = 3
$ws.Cells.Item(240, 6).Value = 3
$ws.Cells.Item(240, 7).Value = "ifop"
$ws.Cells.Item(240, 8).Value = "included"
$ws.Cells.Item(240, 9).Value = "rolling"
$ws.Cells.Item(240, 10).Value = 1117
$ws.Cells.Item(240, 11).Value = 0
$ws.Cells.Item(240, 12).Value = 1
$ws.Cells.Item(240, 13).Value = 0.5
$ws.Cells.Item(240, 14).Value = 0.5
$ws.Cells.Item(240, 15).Value = 0.5
$ws.Cells.Item(240, 16).Value = 11.5
$ws.Cells.Item(240, 17).Value = 4
$ws.Cells.Item(240, 19).Value = 4.5
$ws.Cells.Item(240, 20).Value = 3
$ws.Cells.Item(240, 21).Value = 1
$ws.Cells.Item(240, 22).Value = 28
$ws.Cells.Item(240, 23).Value = 14
$ws.Cells.Item(240, 26).Value = 1.5
$ws.Cells.Item(240, 27).Value = 2
$ws.Cells.Item(240, 28).Value = 17
$ws.Cells.Item(240, 29).Value = 12
$ws.Cells.Item(240, 34).Value = 0.5
$ws.Cells.Item(240, 12).NumberFormat = "0.00"
$ws.Cells.Item(240, 13).NumberFormat = "0.00"

# Row 241 -- poll id 138
$ws.Cells.Item(241, 1).Value = 138
$ws.Cells.Item(241, 2).Value = 2022
$ws.Cells.Item(241, 3).Value = 3
$ws.Cells.Item(241, 4).Value = 2
$ws.Cells.Item(241, 5).Value = 3
$ws.Cells.Item(241, 6).Value = 3
$ws.Cells.Item(241, 7).Value = "ipsos"
$ws.Cells.Item(241, 8).Value = "excluded"
$ws.Cells.Item(241, 9).Value = "regular"
$ws.Cells.Item(241, 10).Value = 2269
$ws.Cells.Item(241, 11).Value = 0
$ws.Cells.Item(241, 12).Value = 1
$ws.Cells.Item(241, 13).Value = 0.5
$ws.Cells.Item(241, 14).Value = 1
$ws.Cells.Item(241, 15).Value = 0.5
$ws.Cells.Item(241, 16).Value = 12
$ws.Cells.Item(241, 17).Value = 4
$ws.Cells.Item(241, 19).Value = 7.5
$ws.Cells.Item(241, 20).Value = 2.5
$ws.Cells.Item(241, 22).Value = 30.5
$ws.Cells.Item(241, 23).Value = 11.5
$ws.Cells.Item(241, 26).Value = 1.5
$ws.Cells.Item(241, 27).Value = 1.5
$ws.Cells.Item(241, 28).Value = 14.5
$ws.Cells.Item(241, 29).Value = 13
$ws.Cells.Item(241, 12).NumberFormat = "0.00"
$ws.Cells.Item(241, 13).NumberFormat = "0.00"

# Row 242 -- poll id 139
$ws.Cells.Item(242, 1).Value = 139
$ws.Cells.Item(242, 2).Value = 2022
$ws.Cells.Item(242, 3).Value = 3
$ws.Cells.Item(242, 4).Value = 2
$ws.Cells.Item(242, 5).Value = 3
$ws.Cells.Item(242, 6).Value = 3
$ws.Cells.Item(242, 7).Value = "opinionway"
$ws.Cells.Item(242, 8).Value = "partially"
$ws.Cells.Item(242, 9).Value = "regular"
$ws.Cells.Item(242, 10).Value = 700
$ws.Cells.Item(242, 11).Value = 0
$ws.Cells.Item(242, 12).Value = 1
$ws.Cells.Item(242, 13).Value = 1
$ws.Cells.Item(242, 14).Value = 1
$ws.Cells.Item(242, 15).Value = "T_1"
$ws.Cells.Item(242, 16).Value = 11
$ws.Cells.Item(242, 17).Value = 4
$ws.Cells.Item(242, 19).Value = 7
$ws.Cells.Item(242, 20).Value = 3
$ws.Cells.Item(242, 22).Value = 28
$ws.Cells.Item(242, 23).Value = 13
$ws.Cells.Item(242, 26).Value = 2
$ws.Cells.Item(242, 27).Value = 2
$ws.Cells.Item(242, 28).Value = 17
$ws.Cells.Item(242, 29).Value = 12
$ws.Cells.Item(242, 12).NumberFormat = "0.00"
$ws.Cells.Item(242, 13).NumberFormat = "0.00"

# Row 243 -- poll id 140
$ws.Cells.Item(243, 1).Value = 140
$ws.Cells.Item(243, 2).Value = 2022
$ws.Cells.Item(243, 3).Value = 3
$ws.Cells.Item(243, 4).Value = 5
$ws.Cells.Item(243, 5).Value = 3
$ws.Cells.Item(243, 6).Value = 6
$ws.Cells.Item(243, 7).Value = "opinionway"
$ws.Cells.Item(243, 8).Value = "partially"
$ws.Cells.Item(243, 9).Value = "regular"
$ws.Cells.Item(243, 10).Value = 700
$ws.Cells.Item(243, 11).Value = 0
$ws.Cells.Item(243, 12).Value = 1
$ws.Cells.Item(243, 13).Value = 1
$ws.Cells.Item(243, 14).Value = 1
$ws.Cells.Item(243, 15).Value = "T_1"
$ws.Cells.Item(243, 16).Value = 10
$ws.Cells.Item(243, 17).Value = 4
$ws.Cells.Item(243, 19).Value = 6
$ws.Cells.Item(243, 20).Value = 3
$ws.Cells.Item(243, 22).Value = 29
$ws.Cells.Item(243, 23).Value = 14
$ws.Cells.Item(243, 26).Value = 2
$ws.Cells.Item(243, 27).Value = 1
$ws.Cells.Item(243, 28).Value = 18
$ws.Cells.Item(243, 29).Value = 12
$ws.Cells.Item(243, 12).NumberFormat = "0.00"
$ws.Cells.Item(243, 13).NumberFormat = "0.00"

# Row 244 -- poll id 141
$ws.Cells.Item(244, 1).Value = 141
$ws.Cells.Item(244, 2).Value = 2022
$ws.Cells.Item(244, 3).Value = 3
$ws.Cells.Item(244, 4).Value = 3
$ws.Cells.Item(244, 5).Value = 3
$ws.Cells.Item(244, 6).Value = 7
$ws.Cells.Item(244, 7).Value = "ifop"
$ws.Cells.Item(244, 8).Value = "included"
$ws.Cells.Item(244, 9).Value = "rolling"
$ws.Cells.Item(244, 10).Value = 1117
$ws.Cells.Item(244, 11).Value = 1
$ws.Cells.Item(244, 12).Value = 1
$ws.Cells.Item(244, 13).Value = 0.5
$ws.Cells.Item(244, 14).Value = 0.5
$ws.Cells.Item(244, 15).Value = "T_0.5"
$ws.Cells.Item(244, 16).Value = 11.5
$ws.Cells.Item(244, 17).Value = 4
$ws.Cells.Item(244, 19).Value = 5
$ws.Cells.Item(244, 20).Value = 2.5
$ws.Cells.Item(244, 22).Value = 30
$ws.Cells.Item(244, 23).Value = 13
$ws.Cells.Item(244, 26).Value = 1.5
$ws.Cells.Item(244, 27).Value = 1.5
$ws.Cells.Item(244, 28).Value = 18
$ws.Cells.Item(244, 29).Value = 12.5
$ws.Cells.Item(244, 12).NumberFormat = "0.00"
$ws.Cells.Item(244, 13).NumberFormat = "0.00"

# Row 245 -- poll id 142
$ws.Cells.Item(245, 1).Value = 142
$ws.Cells.Item(245, 2).Value = 2022
$ws.Cells.Item(245, 3).Value = 3
$ws.Cells.Item(245, 4).Value = 3
$ws.Cells.Item(245, 5).Value = 3
$ws.Cells.Item(245, 6).Value = 7
$ws.Cells.Item(245, 7).Value = "opinionway"
$ws.Cells.Item(245, 8).Value = "partially"
$ws.Cells.Item(245, 9).Value = "rolling"
$ws.Cells.Item(245, 10).Value = 1037
$ws.Cells.Item(245, 11).Value = 1
$ws.Cells.Item(245, 12).Value = 1
$ws.Cells.Item(245, 13).Value = 1
$ws.Cells.Item(245, 14).Value = 1
$ws.Cells.Item(245, 15).Value = 1
$ws.Cells.Item(245, 16).Value = 10
$ws.Cells.Item(245, 17).Value = 5
$ws.Cells.Item(245, 19).Value = 5
$ws.Cells.Item(245, 20).Value = 3
$ws.Cells.Item(245, 22).Value = 30
$ws.Cells.Item(245, 23).Value = 13
$ws.Cells.Item(245, 26).Value = 2
$ws.Cells.Item(245, 27).Value = 1
$ws.Cells.Item(245, 28).Value = 18
$ws.Cells.Item(245, 29).Value = 11
$ws.Cells.Item(245, 12).NumberFormat = "0.00"
$ws.Cells.Item(245, 13).NumberFormat = "0.00"

# Row 246 -- poll id 143
$ws.Cells.Item(246, 1).Value = 143
$ws.Cells.Item(246, 2).Value = 2022
$ws.Cells.Item(246, 3).Value = 3
$ws.Cells.Item(246, 4).Value = 5
$ws.Cells.Item(246, 5).Value = 3
$ws.Cells.Item(246, 6).Value = 6
$ws.Cells.Item(246, 7).Value = "opinionway"
$ws.Cells.Item(246, 8).Value = "partially"
$ws.Cells.Item(246, 9).Value = "regular"
$ws.Cells.Item(246, 10).Value = 700
$ws.Cells.Item(246, 11).Value = 0
$ws.Cells.Item(246, 12).Value = 1
$ws.Cells.Item(246, 13).Value = 1
$ws.Cells.Item(246, 14).Value = 1
$ws.Cells.Item(246, 15).Value = "T_1"
$ws.Cells.Item(246, 16).Value = 11
$ws.Cells.Item(246, 17).Value = 4
$ws.Cells.Item(246, 19).Value = 7
$ws.Cells.Item(246, 20).Value = 3
$ws.Cells.Item(246, 22).Value = 28
$ws.Cells.Item(246, 23).Value = 13
$ws.Cells.Item(246, 26).Value = 2
$ws.Cells.Item(246, 27).Value = 2
$ws.Cells.Item(246, 28).Value = 17
$ws.Cells.Item(246, 29).Value = 12
$ws.Cells.Item(246, 12).NumberFormat = "0.00"
$ws.Cells.Item(246, 13).NumberFormat = "0.00"

# --- Match the workbook's on-screen state after the edit: scrolled so row 224 is the
# first row below the frozen header, with cell A243 selected.
$excel.ActiveWindow.ScrollRow = 224
$ws.Range("A243").Select()
